$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Handout master date placeholder: 4/13/20 -> 4/16/20
# ---------------------------------------------------------------------------
$hm = $p.HandoutMaster
$hf = $hm.HeadersFooters
$dt = $hf.DateAndTime
$dt.UseFormat = 0
$dt.Text = "4/16/20"

# ---------------------------------------------------------------------------
# 2) Slide 14, "Rectangle 5" diagram textbox: swap the IP Header / SRH blocks
#    and update the figure caption text.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$shp14 = $s14.Shapes.Item(4)
$tr14 = $shp14.TextFrame.TextRange

$paraCount = $tr14.Paragraphs().Count
$lines = New-Object 'System.Collections.Generic.List[string]'
for ($i = 1; $i -le $paraCount; $i++) {
    $t = $tr14.Paragraphs($i).Text
    if ($t.Length -gt 0 -and $t[$t.Length - 1] -eq [char]13) {
        $t = $t.Substring(0, $t.Length - 1)
    }
    [void]$lines.Add($t)
}

$newBlock = @(
    "  +---------------------------------------------------------------+",
    "  | IP Header                                                     |",
    "  .  Source IP Address = Endpoint IPv6 Address                    .",
    "  .  Destination IP Address = Sender IPv6 Address                 .",
    "  .  Protocol = UDP                                               .",
    "  .                                                               .",
    "  +---------------------------------------------------------------+",
    "  |  SRH                                                          |",
    "  .  <Segment List>                                               .",
    "  .  END.TSF with Target SID                                      ."
)
for ($i = 0; $i -lt $newBlock.Length; $i++) {
    $lines[$i] = $newBlock[$i]
}

# Figure caption paragraph is the last one.
$lines[$lines.Count - 1] = "    Figure: Probe Message Header for SRv6 with Endpoint Function"

$tr14.Text = [string]::Join("`r", $lines)

# Re-apply bold formatting to "Endpoint" and "Sender" inside the new text.
$paraSource = $tr14.Paragraphs(3)
$paraDest = $tr14.Paragraphs(4)

$idxEndpoint = $paraSource.Text.IndexOf("Endpoint")
$boldEndpoint = $tr14.Characters($paraSource.Start + $idxEndpoint, 8)
$boldEndpoint.Font.Bold = 1

$idxSender = $paraDest.Text.IndexOf("Sender")
$boldSender = $tr14.Characters($paraDest.Start + $idxSender, 6)
$boldSender.Font.Bold = 1

# ---------------------------------------------------------------------------
# 3) Slide 9, "Content Placeholder 2": reduce body text size 18pt -> 16pt
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9 = $shp9.TextFrame.TextRange
$tr9.Font.Size = 16
